$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.620.04'
$ws.Range('E2').Value = '  +0.40%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.961.58'
$ws.Range('E3').Value = '  +1.32%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('E5').Value = '  +1.20%  '

$ws.Range('E6').Value = '  +1.40%  '

$ws.Range('E7').Value = '  +7.57%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  +5.00%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0797'
$ws.Range('E10').Value = '  -5.28%  '

$ws.Range('E11').Value = '  +0.67%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.33'
$ws.Range('E12').Value = '  +7.19%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.838'
$ws.Range('E13').Value = '  +4.48%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.94'
$ws.Range('E14').Value = '  +3.91%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.239.18'
$ws.Range('E15').Value = '  +0.90%  '

$ws.Range('E16').Value = '  +3.80%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.961.67'
$ws.Range('E17').Value = '  +1.79%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.546.45'
$ws.Range('E18').Value = '  +0.52%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.95'
$ws.Range('E19').Value = '  +1.62%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0854'
$ws.Range('E20').Value = '  -0.50%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '230.51'
$ws.Range('E21').Value = '  +1.83%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.09'
$ws.Range('E22').Value = '  +3.11%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.12%  '

$ws.Range('E24').Value = '  +6.15%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  +3.47%  '

$ws.Range('E26').Value = '  +7.50%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.21'
$ws.Range('E27').Value = '  +1.64%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.59'
$ws.Range('E28').Value = '  +0.05%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.43'
$ws.Range('E29').Value = '  +1.53%  '

$ws.Range('E30').Value = '  +11.99%  '

$ws.Range('E31').Value = '  +2.04%  '

$ws.Range('E32').Value = '  +5.97%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0618'
$ws.Range('E33').Value = '  +0.01%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.48'
$ws.Range('E34').Value = '  +8.43%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.55'
$ws.Range('E35').Value = '  +18.46%  '

$ws.Range('E36').Value = '  +6.57%  '

$ws.Range('E37').Value = '  -0.25%  '

$ws.Range('E38').Value = '  -0.08%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.58'
$ws.Range('E39').Value = '  -5.80%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0990'
$ws.Range('E40').Value = '  +2.00%  '

$ws.Range('E41').Value = '  +1.25%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.17'
$ws.Range('E42').Value = '  +2.72%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0211'
$ws.Range('E43').Value = '  +1.81%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.20'
$ws.Range('E44').Value = '  +4.68%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.365.03'
$ws.Range('E45').Value = '  +2.65%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.77'
$ws.Range('E46').Value = '  +4.49%  '

$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.04'
$ws.Range('E47').Value = '  +2.54%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.14'
$ws.Range('E48').Value = '  +1.71%  '

$ws.Range('E49').Value = '  +0.92%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '44.36'
$ws.Range('E50').Value = '  +2.52%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.133.73'
$ws.Range('E51').Value = '  +1.10%  '
